# forests-scraped.xlsx update — 2026-02-03 12:31
# 1) The 7 rows currently sitting in "New" (rows 2-8) are archived to the
#    bottom of "Previously added" (appended as rows 454-460), keeping their
#    original text/hyperlinks.
# 2) "New" rows 2-8 are then overwritten with freshly scraped listings.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# ---- 1) rows moving from "New" -> end of "Previously added" ---------------
$archiveRows = @(
  @{A="https://www.ss.com/msg/lv/real-estate/wood/aizkraukle-and-reg/plavinas/mfgbj.html";        B="6 000 €";   C="Aizkraukle un raj."; D="1 ha.";    E="32420090035"; F=46055.46388888889},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/aluksne-and-reg/jaunaluksnes-pag/inghf.html";    B="33 000 €";  C="Alūksne un raj.";    D="2.30 ha."; E="36560130028"; F=46055.36041666666},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/izvaltas-pag/hmlmn.html";       B="120 000 €"; C="Krāslava un raj.";    D="19 ha.";   E="60640020120"; F=46054.705555555556},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/kuldiga-and-reg/kurmales-pag/chgio.html";        B="13 900 €";  C="Kuldīga un raj.";    D="2 ha.";    E="";            F=46054.82986111111},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/blontu-pag/dioce.html";            B="41 000 €";  C="Ludza un raj.";      D="4 ha.";    E="68440050028"; F=46052.674305555556},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/ogre-and-reg/lielvardes-l-t/lhlxf.html";         B="28 000 €";  C="Ogre un raj.";       D="8 ha.";    E="74840070028"; F=46053.77916666667},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/saldus-and-reg/saldus/kgkjn.html";               B="6 000 €";   C="Saldus un raj.";     D="1.35 ha."; E="84480060140"; F=46055.40763888889}
)

$destStart = 454
for ($i = 0; $i -lt $archiveRows.Count; $i++) {
  $r    = $destStart + $i
  $data = $archiveRows[$i]

  # Column E holds cadastre numbers that are all-digits (e.g. "32420090035");
  # force text so they don't get auto-coerced into numeric cells.
  $ws1.Cells($r,5).NumberFormat = "@"

  $ws1.Cells($r,1).Value = $data.A
  $ws1.Cells($r,2).Value = $data.B
  $ws1.Cells($r,3).Value = $data.C
  $ws1.Cells($r,4).Value = $data.D
  $ws1.Cells($r,5).Value = $data.E
  $ws1.Cells($r,6).Value = $data.F

  $ws1.Hyperlinks.Add($ws1.Cells($r,1), $data.A)

  # Clone formatting from the previous (known-good) row last, so the new
  # row picks up the same style indices (s="3"/"4"/"2") as its neighbours -
  # this also overwrites the generic "Hyperlink" style Hyperlinks.Add just
  # stamped onto column A and the "@" override on column E.
  $ws1.Range("A453:F453").Copy()
  $ws1.Range("A$r`:F$r").PasteSpecial(-4122)
}

# ---- 2) "New" rows 2-8 replaced with the newly scraped listings -----------
$newRows = @(
  @{A="https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/vecsaules-pag/kffcn.html"; B="19 900 €"; C="Bauska un raj.";  D="4 ha.";    E="40920090234"; F=46056.336111111115},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/cesis-and-reg/cesis/cxhdf.html";          B="33 000 €"; C="Cēsis un raj.";    D="1 ha.";    E="42460060005"; F=46056.475694444445},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/upmalas-pag/bhhlnf.html";  B="65 000 €"; C="Preiļi un raj.";   D="14 ha.";   E="76900060029"; F=46055.675},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/dricanu-pag/jhdxd.html";  B="13 000 €"; C="Rēzekne un raj.";  D="4.70 ha."; E="78500020015"; F=46056.495833333334},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/veremu-pag/bgglm.html";   B="20 000 €"; C="Rēzekne un raj.";  D="8.40 ha."; E="78960030042"; F=46055.882638888885},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/valka-and-reg/smiltene/piixf.html";       B="2 000 €";  C="Valka un raj.";    D="0.25 ha."; E="94480020026"; F=46056.513194444444},
  @{A="https://www.ss.com/msg/lv/real-estate/wood/other/hjkix.html";                        B="19 900 €"; C="";                 D="4 ha.";    E="42660060086"; F=46056.30208333333}
)

$newStart = 2
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r    = $newStart + $i
  $data = $newRows[$i]

  # Drop the hyperlink that used to live on this cell before overwriting it.
  $ws2.Cells($r,1).Hyperlinks.Delete()

  # Column E holds cadastre numbers that are all-digits (e.g. "40920090234");
  # force text so they don't get auto-coerced into numeric cells.
  $ws2.Cells($r,5).NumberFormat = "@"

  $ws2.Cells($r,1).Value = $data.A
  $ws2.Cells($r,2).Value = $data.B
  $ws2.Cells($r,3).Value = $data.C
  $ws2.Cells($r,4).Value = $data.D
  $ws2.Cells($r,5).Value = $data.E
  $ws2.Cells($r,6).Value = $data.F

  $ws2.Hyperlinks.Add($ws2.Cells($r,1), $data.A)

  # Re-stamp formatting (keeps s="3"/"4"/"2" consistent) by cloning from the
  # archive sheet's last, known-good row - also cleans up the "Hyperlink"
  # style Hyperlinks.Add just stamped onto column A and the "@" override
  # on column E.
  $ws1.Range("A453:F453").Copy()
  $ws2.Range("A$r`:F$r").PasteSpecial(-4122)
}
